# This script applies the "USCDI5-Sandbox" content update to
# StructureDefinition-us-core-interpreter-required.xlsx:
#   1. Updates the "Description" metadata text (Metadata!B13)
#   2. Updates the three existing "Context" rows (Metadata!B22:B24) and
#      appends a fourth Context row (Metadata!A25:B25) for element:Practitioner
#   3. Shortens the "Definition" text for the root Extension row on the
#      Elements sheet (Elements!M2)
#   4. Updates the Extension.value[x] Type(s) from "code" to "Coding"
#      (Elements!K6), keeping the trailing newline present in the source

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item(1)
$elements = $wb.Worksheets.Item(2)

# 1) Description text (Metadata sheet, row 13)
$metadata.Range("B13").Value = "This extension indicates the individual's need for an interpreter to communicate healthcare information in a language other than the default language of the organization."

# 2) Context values - rewrite the existing 3 rows and add a new 4th row
$metadata.Range("B22").Value = "element:Encounter"
$metadata.Range("B23").Value = "element:Patient"
$metadata.Range("B24").Value = "element:RelatedPerson"

# Duplicate row 24's formatting down into the new row 25 so the new row
# matches the existing Context rows (border/wrap style "s=2")
$metadata.Range("A24:B24").Copy($metadata.Range("A25:B25"))
$metadata.Range("A25").Value = "Context"
$metadata.Range("B25").Value = "element:Practitioner"

# 3) Shorten the Definition text of the "Whether the individual needs an
#    interpreter" row (Elements sheet, row 2, column M)
$elements.Range("M2").Value = "This individual needs an interpreter to communicate healthcare information."

# 4) Extension.value[x] Type(s) changes from "code" to "Coding" (Elements
#    sheet, row 6, column K) - preserve the trailing newline character
$elements.Range("K6").Value = "Coding`n"
